# "Fruta / hortaliza, semanal" — weekly refresh of the Albahaca price sheet.
# A new weekly observation is inserted as row 71 (pushing the previous rows
# 71-78 down to 72-79), and the new row is populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 71; existing rows 71:78 shift to 72:79.
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new weekly record.
$ws.Cells.Item(71, 1).Value  = 1
$ws.Cells.Item(71, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(71, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(71, 4).Value  = 45218
$ws.Cells.Item(71, 5).Value  = 15
$ws.Cells.Item(71, 6).Value  = 100112052
$ws.Cells.Item(71, 7).Value  = "Albahaca"
$ws.Cells.Item(71, 8).Value  = "Sin especificar"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 450
$ws.Cells.Item(71, 11).Value = 1300
$ws.Cells.Item(71, 12).Value = 1500
$ws.Cells.Item(71, 13).Value = 1389
$ws.Cells.Item(71, 14).Value = "$/paquete"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 1389
$ws.Cells.Item(71, 17).Value = 1
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# Keep the date column formatted like the rest of column D.
$ws.Range("D71").NumberFormat = $ws.Range("D72").NumberFormat
